$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B149").Value = 48654
$ws.Range("C149").Value = "CHO-Medimix Sandal with Eladi oils for glowing skin and natural protection Soap-125gms"
$ws.Range("E149").Value = 38.26
$ws.Range("F149").Value = -1
$ws.Range("G149").Value = -32.02

$ws.Range("B150").Value = 63902
$ws.Range("C150").Value = "CHO-Medimix Sandal with Eladi oils for glowing skin and natural protection Soap-125gms"
$ws.Range("E150").Value = 34.04
$ws.Range("F150").Value = 2
$ws.Range("G150").Value = 64.04000000000001

$ws.Range("B161").Value = 53925
$ws.Range("C161").Value = "COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush"
$ws.Range("E161").Value = 79.37
$ws.Range("F161").Value = 1
$ws.Range("G161").Value = 66.44

$ws.Range("B162").Value = 64350
$ws.Range("C162").Value = "COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush"
$ws.Range("E162").Value = 70.63
$ws.Range("F162").Value = 101
$ws.Range("G162").Value = 6710.44

$ws.Range("B163").Value = 57756
$ws.Range("C163").Value = "COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush"
$ws.Range("E163").Value = 79.37
$ws.Range("F163").Value = -100
$ws.Range("G163").Value = -6644

$ws.Range("B264").Value = 64979
$ws.Range("C264").Value = "HIM-BABY CARE GIFT PACK (WW)1"
$ws.Range("E264").Value = 314.41
$ws.Range("F264").Value = 82
$ws.Range("G264").Value = 24251.5

$ws.Range("B265").Value = 48719
$ws.Range("C265").Value = "HIM-BABY CARE GIFT PACK (WW)1"
$ws.Range("E265").Value = 353.35
$ws.Range("F265").Value = -81
$ws.Range("G265").Value = -23955.75

$ws.Range("B279").Value = 64973
$ws.Range("C279").Value = "HIM-GENTLE BABY SOAP 75G"
$ws.Range("E279").Value = 35.4
$ws.Range("F279").Value = 150
$ws.Range("G279").Value = 4995

$ws.Range("B280").Value = 48706
$ws.Range("C280").Value = "HIM-GENTLE BABY SOAP 75G"
$ws.Range("E280").Value = 39.8
$ws.Range("F280").Value = -144
$ws.Range("G280").Value = -4795.2

$ws.Range("B313").Value = 62997
$ws.Range("C313").Value = "HUL-3Roses Dust [C] 500G Relaunch"
$ws.Range("E313").Value = 325.16
$ws.Range("F313").Value = 72
$ws.Range("G313").Value = 22020.48

$ws.Range("B314").Value = 57854
$ws.Range("C314").Value = "HUL-3Roses Dust [C] 500G Relaunch"
$ws.Range("E314").Value = 325.16
$ws.Range("F314").Value = 2
$ws.Range("G314").Value = 611.6799999999999

$ws.Range("B350").Value = 63531
$ws.Range("C350").Value = "HUL-Kissan Pineapple Jam 500G"
$ws.Range("E350").Value = 152.53
$ws.Range("F350").Value = 80
$ws.Range("G350").Value = 11478.4

$ws.Range("B351").Value = 63571
$ws.Range("C351").Value = "HUL-Kissan Pineapple Jam 500G"
$ws.Range("E351").Value = 152.53
$ws.Range("F351").Value = 27
$ws.Range("G351").Value = 3873.96

$ws.Range("B352").Value = 57802
$ws.Range("C352").Value = "HUL-Kissan Pineapple Jam 500G"
$ws.Range("E352").Value = 162.71
$ws.Range("F352").Value = -79
$ws.Range("G352").Value = -11334.92

$ws.Range("B372").Value = 57885
$ws.Range("C372").Value = "HUL-Liril Soap 125 G"
$ws.Range("E372").Value = 62.28
$ws.Range("F372").Value = 4
$ws.Range("G372").Value = 208.52

$ws.Range("B373").Value = 63652
$ws.Range("C373").Value = "HUL-Liril Soap 125 G"
$ws.Range("E373").Value = 55.42
$ws.Range("F373").Value = 250
$ws.Range("G373").Value = 13032.5

$ws.Range("B379").Value = 63564
$ws.Range("C379").Value = "HUL-Lux Radiant Glow 3*150g"
$ws.Range("E379").Value = 137.16
$ws.Range("F379").Value = 57
$ws.Range("G379").Value = 7353.57

$ws.Range("B380").Value = 61608
$ws.Range("C380").Value = "HUL-Lux Radiant Glow 3*150g"
$ws.Range("E380").Value = 154.12
$ws.Range("F380").Value = -56
$ws.Range("G380").Value = -7224.56

$ws.Range("B389").Value = 62865
$ws.Range("C389").Value = "HUL-Rap Refresh Bolt 1Kg"
$ws.Range("E389").Value = 95.34999999999999
$ws.Range("F389").Value = 151
$ws.Range("G389").Value = 12051.31

$ws.Range("B390").Value = 57817
$ws.Range("C390").Value = "HUL-Rap Refresh Bolt 1Kg"
$ws.Range("E390").Value = 95.34999999999999
$ws.Range("F390").Value = 3
$ws.Range("G390").Value = 239.43

$ws.Range("B400").Value = 62933
$ws.Range("C400").Value = "HUL-Sfxl Ew Bale 500G"
$ws.Range("E400").Value = 70.65000000000001
$ws.Range("F400").Value = 146
$ws.Range("G400").Value = 8632.98

$ws.Range("B401").Value = 57835
$ws.Range("C401").Value = "HUL-Sfxl Ew Bale 500G"
$ws.Range("E401").Value = 70.65000000000001
$ws.Range("F401").Value = 1
$ws.Range("G401").Value = 59.13

$ws.Range("B419").Value = 63007
$ws.Range("C419").Value = "HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp"
$ws.Range("E419").Value = 204.69
$ws.Range("F419").Value = 984
$ws.Range("G419").Value = 168588.72

$ws.Range("B420").Value = 57856
$ws.Range("C420").Value = "HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp"
$ws.Range("E420").Value = 204.69
$ws.Range("F420").Value = 2
$ws.Range("G420").Value = 342.66

$ws.Range("B421").Value = 63008
$ws.Range("C421").Value = "HUL-Surf Exl Mtc Liq Tl 1 Ltr Cp"
$ws.Range("E421").Value = 180.62
$ws.Range("F421").Value = 504
$ws.Range("G421").Value = 76189.67999999999

$ws.Range("B422").Value = 57857
$ws.Range("C422").Value = "HUL-Surf Exl Mtc Liq Tl 1 Ltr Cp"
$ws.Range("E422").Value = 180.62
$ws.Range("F422").Value = 3
$ws.Range("G422").Value = 453.51

$ws.Range("B431").Value = 53082
$ws.Range("C431").Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Range("E431").Value = 71.05
$ws.Range("F431").Value = 1
$ws.Range("G431").Value = 59.47

$ws.Range("B432").Value = 63102
$ws.Range("C432").Value = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Range("E432").Value = 71.05
$ws.Range("F432").Value = 36
$ws.Range("G432").Value = 2140.92

$ws.Range("B457").Value = 63681
$ws.Range("C457").Value = "JLM-MBD Shiny Toothbrush Safari"
$ws.Range("E457").Value = 23.84
$ws.Range("F457").Value = 65
$ws.Range("G457").Value = 1457.3

$ws.Range("B458").Value = 31930
$ws.Range("C458").Value = "JLM-MBD Shiny Toothbrush Safari"
$ws.Range("E458").Value = 26.8
$ws.Range("F458").Value = -62
$ws.Range("G458").Value = -1390.04

$ws.Range("B579").Value = 65069
$ws.Range("C579").Value = "CRE-Bourbon 100gm"
$ws.Range("E579").Value = 14.3
$ws.Range("F579").Value = 172
$ws.Range("G579").Value = 2313.4

$ws.Range("B580").Value = 53757
$ws.Range("C580").Value = "CRE-Bourbon 100gm"
$ws.Range("E580").Value = 16.08
$ws.Range("F580").Value = -159
$ws.Range("G580").Value = -2138.55

$ws.Range("B590").Value = 64922
$ws.Range("C590").Value = "CRE-Cremica Golden Bytes Rich Butter 200Gm"
$ws.Range("E590").Value = 20.98
$ws.Range("F590").Value = 207
$ws.Range("G590").Value = 4084.11

$ws.Range("B591").Value = 45706
$ws.Range("C591").Value = "CRE-Cremica Golden Bytes Rich Butter 200Gm"
$ws.Range("E591").Value = 23.58
$ws.Range("F591").Value = -202
$ws.Range("G591").Value = -3985.46

$ws.Range("B593").Value = 45718
$ws.Range("C593").Value = "CRE-Cremica Honey Oatmeal Cookies 50 +25 Gm"
$ws.Range("E593").Value = 19.38
$ws.Range("F593").Value = -294
$ws.Range("G593").Value = -4768.68

$ws.Range("B594").Value = 64927
$ws.Range("C594").Value = "CRE-Cremica Honey Oatmeal Cookies 50 +25 Gm"
$ws.Range("E594").Value = 17.26
$ws.Range("F594").Value = 295
$ws.Range("G594").Value = 4784.9

$ws.Range("B601").Value = 64919
$ws.Range("C601").Value = "CRE-Cremica Pista Almond Cookies (75 +25Gm)"
$ws.Range("E601").Value = 27.97
$ws.Range("F601").Value = 224
$ws.Range("G601").Value = 5891.2

$ws.Range("B602").Value = 45702
$ws.Range("C602").Value = "CRE-Cremica Pista Almond Cookies (75 +25Gm)"
$ws.Range("E602").Value = 31.43
$ws.Range("F602").Value = -215
$ws.Range("G602").Value = -5654.5

$ws.Range("B709").Value = 64833
$ws.Range("C709").Value = "Rasna 32 Glass Shikanji Nimbupani"
$ws.Range("E709").Value = 34.9
$ws.Range("F709").Value = 99
$ws.Range("G709").Value = 3250.17

$ws.Range("B710").Value = 60025
$ws.Range("C710").Value = "Rasna 32 Glass Shikanji Nimbupani"
$ws.Range("E710").Value = 37.22
$ws.Range("F710").Value = -98
$ws.Range("G710").Value = -3217.34

$ws.Range("B720").Value = 64830
$ws.Range("C720").Value = "Rasna Nagpur Orange (32 Glass)"
$ws.Range("E720").Value = 34.9
$ws.Range("F720").Value = 117
$ws.Range("G720").Value = 3841.11

$ws.Range("B721").Value = 60022
$ws.Range("C721").Value = "Rasna Nagpur Orange (32 Glass)"
$ws.Range("E721").Value = 37.22
$ws.Range("F721").Value = -113
$ws.Range("G721").Value = -3709.79

$ws.Range("B872").Value = 65362
$ws.Range("C872").Value = "Shankys Tip Top Hing Jeera Peanut/ Salted Peanut 200 Gm"
$ws.Range("E872").Value = 43.44
$ws.Range("F872").Value = 2
$ws.Range("G872").Value = 81.73999999999999

$ws.Range("B873").Value = 65079
$ws.Range("C873").Value = "Shankys Tip Top Hing Jeera Peanut/ Salted Peanut 200 Gm"
$ws.Range("E873").Value = 43.44
$ws.Range("F873").Value = 21
$ws.Range("G873").Value = 858.27
